$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.073.26'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '3.140.16'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'589.90"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = "'146.58"
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.129.78'
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = '  +5.54%  '
$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = "'37.39"
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '3.665.61'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = "'7.25"
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("D18").Value = '63.891.58'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").Value = '3.136.31'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").Value = "'466.71"
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = "'14.38"
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").Value = "'0.730"
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").Value = "'7.55"
$ws.Range("E23").Value = '  +1.38%  '
$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = '  +12.83%  '
$ws.Range("D25").Value = "'13.08"
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("D26").Value = "'80.90"
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("D28").Value = "'9.97"
$ws.Range("E28").Value = '  +10.71%  '
$ws.Range("D29").Value = "'2.70"
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("E33").Value = '  +4.64%  '
$ws.Range("D34").Value = "'27.62"
$ws.Range("E34").Value = '  +3.76%  '
$ws.Range("D35").Value = '0.0₃0851'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = "'6.15"
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = "'3.27"
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("D40").Value = "'458.54"
$ws.Range("E40").Value = '  +5.48%  '
$ws.Range("E43").Value = '  +7.47%  '
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("D45").Value = '2.886.05'
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("D46").Value = "'39.89"
$ws.Range("E46").Value = '  +11.44%  '
$ws.Range("D48").Value = "'133.76"
$ws.Range("E48").Value = '  +8.49%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  +3.44%  '

# Rows 30-32: rotate NEARProtocol/ImmutableX/FirstDigitalUSD with updated prices
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = "'7.27"
$ws.Range("E31").Value = '  +6.35%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = '  +0.61%  '

# Rows 41-42: swap OKB/Cosmos with updated prices
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = "'9.33"
$ws.Range("E41").Value = '  +6.95%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'51.19"
$ws.Range("E42").Value = '  +1.90%  '
